$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H2").Value = 129.83333
$ws_ALC.Range("I2").Value = 129.83333
$ws_ALC.Range("K2").Value = 129.83333
$ws_ALC.Range("M2").Value = -16.83332999999999

$ws_ALC.Range("H15").Value = 241096.45
$ws_ALC.Range("I15").Value = 241096.45
$ws_ALC.Range("K15").Value = 723289.3500000001
$ws_ALC.Range("M15").Value = -723120.3500000001

$ws_ALC.Range("H32").Value = 2556.7144
$ws_ALC.Range("J32").Value = 2833.1667
$ws_ALC.Range("L32").Value = 2833.1667
$ws_ALC.Range("N32").Value = -3485.1667

$ws_ALC.Range("H106").Value = 22222888
$ws_ALC.Range("I106").Value = 33333832
$ws_ALC.Range("J106").Value = 1000
$ws_ALC.Range("K106").Value = 33333832
$ws_ALC.Range("L106").Value = 1000
$ws_ALC.Range("M106").Value = -33333201
$ws_ALC.Range("N106").Value = -2262

$ws_ALC.Range("H107").Value = 2019.2632
$ws_ALC.Range("I107").Value = 2021.5294
$ws_ALC.Range("K107").Value = 2021.5294
$ws_ALC.Range("M107").Value = -101.5293999999999

$ws_ALC.Range("H112").Value = 4279.0684
$ws_ALC.Range("J112").Value = 4854.972
$ws_ALC.Range("L112").Value = 14564.916
$ws_ALC.Range("N112").Value = -16780.916

$ws_ALC.Range("H141").Value = 3804.0908
$ws_ALC.Range("I141").Value = 3986
$ws_ALC.Range("J141").Value = 1985
$ws_ALC.Range("K141").Value = 11958
$ws_ALC.Range("L141").Value = 5955
$ws_ALC.Range("M141").Value = -6778
$ws_ALC.Range("N141").Value = -16315

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 2390.3823
$ws_ARM.Range("I2").Value = 2337.6897
$ws_ARM.Range("K2").Value = 2337.6897
$ws_ARM.Range("M2").Value = -2224.6897

$ws_ARM.Range("H32").Value = 8112.65
$ws_ARM.Range("I32").Value = 7016.6045
$ws_ARM.Range("K32").Value = 7016.6045
$ws_ARM.Range("M32").Value = -6729.6045

$ws_ARM.Range("H92").Value = 55742
$ws_ARM.Range("J92").Value = 55742
$ws_ARM.Range("L92").Value = 55742
$ws_ARM.Range("N92").Value = -60734

$ws_ARM.Range("H116").Value = 2390.3823
$ws_ARM.Range("I116").Value = 2337.6897
$ws_ARM.Range("K116").Value = 2337.6897
$ws_ARM.Range("M116").Value = -43.6896999999999

$ws_ARM.Range("H132").Value = 192493.73
$ws_ARM.Range("I132").Value = 306431.72
$ws_ARM.Range("J132").Value = 4496.05
$ws_ARM.Range("K132").Value = 919295.1599999999
$ws_ARM.Range("L132").Value = 13488.15
$ws_ARM.Range("M132").Value = -916765.1599999999
$ws_ARM.Range("N132").Value = -18548.15

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 2390.3823
$ws_BSM.Range("I3").Value = 2337.6897
$ws_BSM.Range("K3").Value = 2337.6897
$ws_BSM.Range("M3").Value = -2223.6897

$ws_BSM.Range("H94").Value = 2249.8708
$ws_BSM.Range("I94").Value = 2373.5
$ws_BSM.Range("J94").Value = 1096
$ws_BSM.Range("K94").Value = 2373.5
$ws_BSM.Range("L94").Value = 1096
$ws_BSM.Range("M94").Value = -1922.5
$ws_BSM.Range("N94").Value = -1998

$ws_BSM.Range("H132").Value = 137983
$ws_BSM.Range("J132").Value = 137983
$ws_BSM.Range("L132").Value = 137983
$ws_BSM.Range("N132").Value = -148103

$ws_BSM.Range("H134").Value = 10381037
$ws_BSM.Range("I134").Value = 2859173
$ws_BSM.Range("K134").Value = 8577519
$ws_BSM.Range("M134").Value = -8574984

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H5").Value = 1826
$ws_CRP.Range("I5").Value = 435
$ws_CRP.Range("J5").Value = 5999
$ws_CRP.Range("K5").Value = 435
$ws_CRP.Range("L5").Value = 5999
$ws_CRP.Range("M5").Value = -323
$ws_CRP.Range("N5").Value = -6223

$ws_CRP.Range("H31").Value = 6780.9375
$ws_CRP.Range("I31").Value = 2473.3137
$ws_CRP.Range("J31").Value = 14356.414
$ws_CRP.Range("K31").Value = 2473.3137
$ws_CRP.Range("L31").Value = 14356.414
$ws_CRP.Range("M31").Value = -2178.3137
$ws_CRP.Range("N31").Value = -14946.414

$ws_CRP.Range("H34").Value = 6780.9375
$ws_CRP.Range("I34").Value = 2473.3137
$ws_CRP.Range("J34").Value = 14356.414
$ws_CRP.Range("K34").Value = 2473.3137
$ws_CRP.Range("L34").Value = 14356.414
$ws_CRP.Range("M34").Value = -2271.3137
$ws_CRP.Range("N34").Value = -14760.414

$ws_CRP.Range("H58").Value = 3363.5151
$ws_CRP.Range("I58").Value = 3292.36
$ws_CRP.Range("J58").Value = 3585.875
$ws_CRP.Range("K58").Value = 3292.36
$ws_CRP.Range("L58").Value = 3585.875
$ws_CRP.Range("M58").Value = -3089.36
$ws_CRP.Range("N58").Value = -3991.875

$ws_CRP.Range("H92").Value = 129995
$ws_CRP.Range("J92").Value = 129995
$ws_CRP.Range("L92").Value = 129995
$ws_CRP.Range("N92").Value = -134987

$ws_CRP.Range("H94").Value = 1056.3334
$ws_CRP.Range("I94").Value = 904.6
$ws_CRP.Range("K94").Value = 904.6
$ws_CRP.Range("M94").Value = -453.6

$ws_CRP.Range("H107").Value = 1029.2222
$ws_CRP.Range("I107").Value = 793.7619
$ws_CRP.Range("J107").Value = 1853.3334
$ws_CRP.Range("K107").Value = 793.7619
$ws_CRP.Range("L107").Value = 1853.3334
$ws_CRP.Range("M107").Value = 1126.2381
$ws_CRP.Range("N107").Value = -5693.3334

$ws_CRP.Range("H132").Value = 3042.093
$ws_CRP.Range("I132").Value = 3061.4546
$ws_CRP.Range("J132").Value = 2978.2
$ws_CRP.Range("K132").Value = 9184.363799999999
$ws_CRP.Range("L132").Value = 8934.599999999999
$ws_CRP.Range("M132").Value = -6654.363799999999
$ws_CRP.Range("N132").Value = -13994.6

$ws_CRP.Range("H134").Value = 2594.6758
$ws_CRP.Range("I134").Value = 2034.5769
$ws_CRP.Range("K134").Value = 6103.7307
$ws_CRP.Range("M134").Value = -3568.7307

$ws_CRP.Range("H136").Value = 3363.5151
$ws_CRP.Range("I136").Value = 3292.36
$ws_CRP.Range("J136").Value = 3585.875
$ws_CRP.Range("K136").Value = 9877.08
$ws_CRP.Range("L136").Value = 10757.625
$ws_CRP.Range("M136").Value = -7327.08
$ws_CRP.Range("N136").Value = -15857.625

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 948.3143
$ws_CUL.Range("J5").Value = 1544.3636
$ws_CUL.Range("L5").Value = 4633.0908
$ws_CUL.Range("N5").Value = -4857.0908

$ws_CUL.Range("H12").Value = 526349.8
$ws_CUL.Range("I12").Value = 28.833334
$ws_CUL.Range("J12").Value = 769267.25
$ws_CUL.Range("K12").Value = 86.50000199999999
$ws_CUL.Range("L12").Value = 2307801.75
$ws_CUL.Range("M12").Value = 86.49999800000001
$ws_CUL.Range("N12").Value = -2308147.75

$ws_CUL.Range("H131").Value = 2555.6765
$ws_CUL.Range("J131").Value = 2500.8276
$ws_CUL.Range("L131").Value = 7502.4828
$ws_CUL.Range("N131").Value = -17582.4828

$ws_CUL.Range("H135").Value = 948.3143
$ws_CUL.Range("J135").Value = 1544.3636
$ws_CUL.Range("L135").Value = 13899.2724
$ws_CUL.Range("N135").Value = -18969.2724

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H70").Value = 4493.0977
$ws_GSM.Range("J70").Value = 4493.4
$ws_GSM.Range("L70").Value = 4493.4
$ws_GSM.Range("N70").Value = -5033.4

$ws_GSM.Range("H73").Value = 4493.0977
$ws_GSM.Range("J73").Value = 4493.4
$ws_GSM.Range("L73").Value = 4493.4
$ws_GSM.Range("N73").Value = -6365.4

$ws_GSM.Range("H102").Value = 1722.2941
$ws_GSM.Range("I102").Value = 1671.4546
$ws_GSM.Range("K102").Value = 1671.4546
$ws_GSM.Range("M102").Value = -49.45460000000003

$ws_GSM.Range("H108").Value = 54999.75
$ws_GSM.Range("J108").Value = 54999.668
$ws_GSM.Range("L108").Value = 54999.668
$ws_GSM.Range("N108").Value = -62679.668

$ws_GSM.Range("H113").Value = 2765.5557
$ws_GSM.Range("I113").Value = 2751.4
$ws_GSM.Range("K113").Value = 2751.4
$ws_GSM.Range("M113").Value = -581.4000000000001

$ws_GSM.Range("H139").Value = 137781.6
$ws_GSM.Range("J139").Value = 137781.6
$ws_GSM.Range("L139").Value = 137781.6
$ws_GSM.Range("N139").Value = -148061.6

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H46").Value = 2311.0557
$ws_LTW.Range("I46").Value = 959.6
$ws_LTW.Range("J46").Value = 4000.375
$ws_LTW.Range("K46").Value = 959.6
$ws_LTW.Range("L46").Value = 4000.375
$ws_LTW.Range("M46").Value = -771.6
$ws_LTW.Range("N46").Value = -4376.375

$ws_LTW.Range("H61").Value = 4187.1333
$ws_LTW.Range("I61").Value = 3292.8462
$ws_LTW.Range("K61").Value = 3292.8462
$ws_LTW.Range("M61").Value = -3090.8462

$ws_LTW.Range("H82").Value = 3397.7368
$ws_LTW.Range("I82").Value = 3624.4285
$ws_LTW.Range("K82").Value = 3624.4285
$ws_LTW.Range("M82").Value = -3263.4285

$ws_LTW.Range("H85").Value = 3397.7368
$ws_LTW.Range("I85").Value = 3624.4285
$ws_LTW.Range("K85").Value = 3624.4285
$ws_LTW.Range("M85").Value = -2376.4285

$ws_LTW.Range("H94").Value = 50161.5
$ws_LTW.Range("J94").Value = 50161.5
$ws_LTW.Range("L94").Value = 50161.5
$ws_LTW.Range("N94").Value = -51513.5

$ws_LTW.Range("H100").Value = 0
$ws_LTW.Range("J100").Value = 0
$ws_LTW.Range("L100").Value = 0
$ws_LTW.Range("N100").ClearContents()

$ws_LTW.Range("H104").Value = 63499.5
$ws_LTW.Range("J104").Value = 63499.5
$ws_LTW.Range("L104").Value = 63499.5
$ws_LTW.Range("N104").Value = -70487.5

$ws_LTW.Range("H113").Value = 4187.1333
$ws_LTW.Range("I113").Value = 3292.8462
$ws_LTW.Range("K113").Value = 3292.8462
$ws_LTW.Range("M113").Value = -1122.8462

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H81").Value = 44759.58
$ws_WVR.Range("I81").Value = 63173.47
$ws_WVR.Range("K81").Value = 126346.94
$ws_WVR.Range("M81").Value = -125285.94

$ws_WVR.Range("H84").Value = 44759.58
$ws_WVR.Range("I84").Value = 63173.47
$ws_WVR.Range("K84").Value = 631734.7
$ws_WVR.Range("M84").Value = -626430.7

$ws_WVR.Range("H109").Value = 101666.336
$ws_WVR.Range("J109").Value = 101666.336
$ws_WVR.Range("L109").Value = 101666.336
$ws_WVR.Range("N109").Value = -104440.336

$ws_WVR.Range("H113").Value = 450.16666
$ws_WVR.Range("I113").Value = 450.2
$ws_WVR.Range("K113").Value = 1350.6
$ws_WVR.Range("M113").Value = 819.4000000000001

$ws_WVR.Range("H124").Value = 53732.75
$ws_WVR.Range("J124").Value = 53732.75
$ws_WVR.Range("L124").Value = 53732.75
$ws_WVR.Range("N124").Value = -63552.75
